# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column D describes "municipio-nombre" metadata; it moves from being
# modelled as an iaest-measure (medida / xsd:int) to a sdmx dimension
# (dim) mapped through a URI-Municipio column.
#
# Column I describes "estado-del-edificio" metadata; it moves from being
# modelled as an iaest-dimension (dim / skos:Concept, with an external
# mapping-estado-del-edificio.xlsx mapping file) to an iaest-measure
# (medida / xsd:int), so the mapping file reference in I5 is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

$ws.Range("I2").Value = "iaest-measure:estado-del-edificio"
$ws.Range("I3").Value = "medida"
$ws.Range("I4").Value = "xsd:int"

$ws.Range("I5").Clear()
